# Update cryptos price (D) and volume-change (E) columns to refreshed snapshot values.
# Values are written as text (matching the original inline-string cells) by temporarily
# forcing a Text number format, then restoring the default ("Normal") style so no extra
# formatting is left behind on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue "D2" "50.912.47"
Set-TextValue "E2" "  -1.28%  "
Set-TextValue "D3" "2.938.76"
Set-TextValue "E3" "  -1.42%  "
Set-TextValue "D4" "1.00"
Set-TextValue "E4" "  -0.01%  "
Set-TextValue "D5" "375.84"
Set-TextValue "E5" "  -1.78%  "
Set-TextValue "D6" "101.08"
Set-TextValue "E6" "  -2.36%  "
Set-TextValue "E7" "  -1.38%  "
Set-TextValue "E8" "  -0.07%  "
Set-TextValue "D9" "0.588"
Set-TextValue "E9" "  -0.91%  "
Set-TextValue "D10" "36.25"
Set-TextValue "E10" "  -2.26%  "
Set-TextValue "E11" "  -0.63%  "
Set-TextValue "D12" "0.0848"
Set-TextValue "E12" "  -0.02%  "
Set-TextValue "D13" "3.394.62"
Set-TextValue "E13" "  -1.64%  "
Set-TextValue "D14" "18.10"
Set-TextValue "E14" "  -0.92%  "
Set-TextValue "D15" "7.58"
Set-TextValue "E15" "  +0.16%  "
Set-TextValue "E16" "  +52.20%  "
Set-TextValue "D17" "2.938.88"
Set-TextValue "E17" "  -1.38%  "
Set-TextValue "D18" "0.996"
Set-TextValue "E18" "  -0.56%  "
Set-TextValue "D19" "50.882.25"
Set-TextValue "E19" "  -1.24%  "
Set-TextValue "D20" "3.07"
Set-TextValue "E20" "  -5.87%  "
Set-TextValue "D21" "12.47"
Set-TextValue "E21" "  -2.64%  "
Set-TextValue "D22" "0.0₃0955"
Set-TextValue "E22" "  -0.89%  "
Set-TextValue "D23" "266.23"
Set-TextValue "E23" "  +1.57%  "
Set-TextValue "D24" "69.05"
Set-TextValue "E24" "  -0.17%  "
Set-TextValue "D25" "3.19"
Set-TextValue "E25" "  +9.53%  "
Set-TextValue "D26" "8.12"
Set-TextValue "E26" "  -1.28%  "
Set-TextValue "D27" "7.40"
Set-TextValue "E27" "  -2.89%  "
Set-TextValue "D28" "0.999"
Set-TextValue "D29" "25.65"
Set-TextValue "E29" "  -1.43%  "
Set-TextValue "D30" "0.162"
Set-TextValue "E30" "  -4.48%  "
Set-TextValue "E31" "  -7.46%  "
Set-TextValue "D32" "9.97"
Set-TextValue "E32" "  +0.97%  "
Set-TextValue "D33" "50.98"
Set-TextValue "E33" "  -0.03%  "
Set-TextValue "D34" "2.05"
Set-TextValue "E34" "  -0.68%  "
Set-TextValue "D35" "33.35"
Set-TextValue "E35" "  -3.44%  "
Set-TextValue "D36" "0.0441"
Set-TextValue "E36" "  -2.56%  "
Set-TextValue "E37" "  +0.07%  "
Set-TextValue "E38" "  +4.40%  "
Set-TextValue "D39" "0.115"
Set-TextValue "E39" "  -0.42%  "
Set-TextValue "D40" "16.46"
Set-TextValue "E40" "  -2.77%  "
Set-TextValue "E41" "  -1.19%  "
Set-TextValue "D42" "2.46"
Set-TextValue "E42" "  -4.48%  "
Set-TextValue "D43" "119.90"
Set-TextValue "E43" "  -2.34%  "
Set-TextValue "D44" "21.15"
Set-TextValue "E44" "  -2.23%  "
Set-TextValue "D45" "3.40"
Set-TextValue "E45" "  +2.92%  "
Set-TextValue "E46" "  -0.47%  "
Set-TextValue "E47" "  -1.35%  "
Set-TextValue "D48" "2.33"
Set-TextValue "E48" "  -1.64%  "
Set-TextValue "D49" "1.990.37"
Set-TextValue "E49" "  -2.10%  "
Set-TextValue "D50" "0.0327"
Set-TextValue "E50" "  -1.67%  "
Set-TextValue "D51" "5.20"
Set-TextValue "E51" "  +1.09%  "
